# Auto-generated-then-reviewed PowerShell Excel COM-interop script
# Applies the weekly NYPD CompStat data refresh described by the commit diff:
#  - bumps the report Volume/Number and the covered week date range
#  - widens column H to match column E
#  - rewrites the Crime Complaints data block (rows 15-33) with the new weekly figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume number + reporting week dates (rich-text cells in A8 / C9) ---
$ws.Range("A8").Value = "Volume 32   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# --- Column H width now matches column E (both bestFit to the wider value) ---
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Crime Complaints data block (rows 15-33) ---

# Row 15
$ws.Range("A15").Value = "Rape"
$ws.Range("C15").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 2.0
$ws.Range("G15").Value = 2.0
$ws.Range("H15").Value = 0.0
$ws.Range("I15").Value = 2.0
$ws.Range("J15").Value = 3.0
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = "***.*"
$ws.Range("M15").Value = "***.*"
$ws.Range("N15").Value = -33.333333333333

# Row 16
$ws.Range("A16").Value = "Robbery"
$ws.Range("C16").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 4.0
$ws.Range("G16").Value = 4.0
$ws.Range("H16").Value = 0.0
$ws.Range("I16").Value = 8.0
$ws.Range("J16").Value = 11.0
$ws.Range("K16").Value = -27.272727272727
$ws.Range("L16").Value = -42.857142857142
$ws.Range("M16").Value = -27.272727272727
$ws.Range("N16").Value = -88.732394366197

# Row 17
$ws.Range("A17").Value = "Fel. Assault"
$ws.Range("C17").Value = 4.0
$ws.Range("D17").Value = 4.0
$ws.Range("F15").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = 0.0
$ws.Range("H15").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 10.0
$ws.Range("G17").Value = 6.0
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 15.0
$ws.Range("J17").Value = 7.0
$ws.Range("K17").Value = 114.285714285714
$ws.Range("L17").Value = 25.0
$ws.Range("M17").Value = -11.764705882352
$ws.Range("N17").Value = -37.5

# Row 18
$ws.Range("A18").Value = "Burglary"
$ws.Range("C18").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 4.0
$ws.Range("G18").Value = 5.0
$ws.Range("H18").Value = -20.0
$ws.Range("I18").Value = 6.0
$ws.Range("J18").Value = 11.0
$ws.Range("K18").Value = -45.454545454545
$ws.Range("L18").Value = -57.142857142857
$ws.Range("M18").Value = -57.142857142857
$ws.Range("N18").Value = -85.365853658536

# Row 19
$ws.Range("A19").Value = "Gr. Larceny"
$ws.Range("C19").Value = 9.0
$ws.Range("D19").Value = 13.0
$ws.Range("E19").Value = -30.76923076923
$ws.Range("F19").Value = 41.0
$ws.Range("G19").Value = 45.0
$ws.Range("H19").Value = -8.888888888888
$ws.Range("I19").Value = 62.0
$ws.Range("J19").Value = 63.0
$ws.Range("K19").Value = -1.587301587301
$ws.Range("L19").Value = -8.823529411764
$ws.Range("M19").Value = -18.421052631578
$ws.Range("N19").Value = -30.337078651685

# Row 20
$ws.Range("A20").Value = "G.L.A."
$ws.Range("C20").Value = "0"
$ws.Range("D20").Value = 1.0
$ws.Range("E20").Value = -100.0
$ws.Range("F20").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("G20").Value = 7.0
$ws.Range("H20").Value = -100.0
$ws.Range("I20").Value = 2.0
$ws.Range("J20").Value = 9.0
$ws.Range("K20").Value = -77.777777777777
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = 0.0
$ws.Range("N20").Value = -96.428571428571

# Row 21
$ws.Range("A21").Value = "TOTAL"
$ws.Range("C21").Value = 15.0
$ws.Range("D21").Value = 18.0
$ws.Range("E21").Value = -16.666666666666
$ws.Range("F21").Value = 61.0
$ws.Range("G21").Value = 69.0
$ws.Range("H21").Value = -11.59420289855
$ws.Range("I21").Value = 95.0
$ws.Range("J21").Value = 104.0
$ws.Range("K21").Value = -8.653846153846
$ws.Range("L21").Value = -14.414414414414
$ws.Range("M21").Value = -20.833333333333
$ws.Range("N21").Value = -66.549295774647

# Row 22
$ws.Range("A22").Value = "Transit"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1.0
$ws.Range("G22").Value = 2.0
$ws.Range("H22").Value = -50.0
$ws.Range("I22").Value = 2.0
$ws.Range("J22").Value = 2.0
$ws.Range("K22").Value = 0.0
$ws.Range("L22").Value = -60.0
$ws.Range("M22").Value = 100.0
$ws.Range("N22").Value = "***.*"

# Row 23
$ws.Range("A23").Value = "Housing"
$ws.Range("C23").Value = 4.0
$ws.Range("F15").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 2.0
$ws.Range("E23").Value = 100.0
$ws.Range("F23").Value = 6.0
$ws.Range("G23").Value = 5.0
$ws.Range("H23").Value = 20.0
$ws.Range("I23").Value = 8.0
$ws.Range("J23").Value = 6.0
$ws.Range("K23").Value = 33.333333333333
$ws.Range("L23").Value = 166.666666666667
$ws.Range("M23").Value = 166.666666666667
$ws.Range("N23").Value = "***.*"

# Row 24
$ws.Range("A24").Value = "Petit Larceny"
$ws.Range("C24").Value = 16.0
$ws.Range("D24").Value = 10.0
$ws.Range("E24").Value = 60.0
$ws.Range("F24").Value = 67.0
$ws.Range("G24").Value = 44.0
$ws.Range("H24").Value = 52.272727272727
$ws.Range("I24").Value = 91.0
$ws.Range("J24").Value = 67.0
$ws.Range("K24").Value = 35.820895522388
$ws.Range("L24").Value = 44.444444444444
$ws.Range("M24").Value = 7.058823529411
$ws.Range("N24").Value = "***.*"

# Row 25
$ws.Range("A25").Value = "Retail Theft"
$ws.Range("C25").Value = 3.0
$ws.Range("D25").Value = 4.0
$ws.Range("E25").Value = -25.0
$ws.Range("F25").Value = 24.0
$ws.Range("G25").Value = 18.0
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 35.0
$ws.Range("J25").Value = 29.0
$ws.Range("K25").Value = 20.689655172413
$ws.Range("L25").Value = -7.894736842105
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"

# Row 26
$ws.Range("A26").Value = "Misd. Assault"
$ws.Range("C26").Value = 7.0
$ws.Range("D26").Value = 6.0
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 18.0
$ws.Range("G26").Value = 29.0
$ws.Range("H26").Value = -37.931034482758
$ws.Range("I26").Value = 26.0
$ws.Range("J26").Value = 42.0
$ws.Range("K26").Value = -38.095238095238
$ws.Range("L26").Value = -35.0
$ws.Range("M26").Value = -21.212121212121
$ws.Range("N26").Value = "***.*"

# Row 27
$ws.Range("A27").Value = "UCR Rape*"
$ws.Range("C27").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 2.0
$ws.Range("G27").Value = 2.0
$ws.Range("H27").Value = 0.0
$ws.Range("I27").Value = 2.0
$ws.Range("J27").Value = 3.0
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = "***.*"
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28
$ws.Range("A28").Value = "Other Sex Crimes"
$ws.Range("C28").Value = 1.0
$ws.Range("D28").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 4.0
$ws.Range("G28").Value = 5.0
$ws.Range("H28").Value = -20.0
$ws.Range("I28").Value = 4.0
$ws.Range("J28").Value = 6.0
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"

# Row 29
$ws.Range("A29").Value = "Shooting Vic."
$ws.Range("C29").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = "0"
$ws.Range("H29").Value = "***.*"
$ws.Range("I29").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("J29").Value = "0"
$ws.Range("K29").Value = "***.*"
$ws.Range("L29").Value = "***.*"
$ws.Range("M29").Value = "***.*"
$ws.Range("N29").Value = -50.0

# Row 30
$ws.Range("A30").Value = "Shooting Inc."
$ws.Range("C30").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("G30").Value = "0"
$ws.Range("H30").Value = "***.*"
$ws.Range("I30").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("J30").Value = "0"
$ws.Range("K30").Value = "***.*"
$ws.Range("L30").Value = "***.*"
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = -50.0

# Row 31
$ws.Range("A31").Value = "Hate Crimes"
$ws.Range("C31").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D31").Value = "0"
$ws.Range("E31").Value = "***.*"
$ws.Range("F31").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("G31").Value = "0"
$ws.Range("H31").Value = "***.*"
$ws.Range("I31").Value = 1.0
$ws.Range("F15").Copy()
$ws.Range("I31").PasteSpecial(-4122)
$ws.Range("J31").Value = "0"
$ws.Range("K31").Value = "***.*"
$ws.Range("L31").Value = -50.0
$ws.Range("H15").Copy()
$ws.Range("L31").PasteSpecial(-4122)
$ws.Range("M31").Value = "***.*"
$ws.Range("N31").Value = "***.*"

# Row 33
$ws.Range("L33").Value = -100.0
$ws.Range("H15").Copy()
$ws.Range("L33").PasteSpecial(-4122)
